$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells that would otherwise be
# auto-parsed as numbers (losing trailing zeros / exact text), then
# restore default (unstyled) cell style so no stray style index is added.
$priceCells = @("D2","D3","D4","D5","D6","D8","D9","D10","D12","D13","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '59.901.42'
$ws.Range("E2").Value = '  -6.05%  '

$ws.Range("D3").Value = '3.344.89'
$ws.Range("E3").Value = '  -2.35%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '563.35'
$ws.Range("E5").Value = '  -3.07%  '

$ws.Range("D6").Value = '130.58'
$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '3.349.97'
$ws.Range("E8").Value = '  -2.19%  '

$ws.Range("D9").Value = '0.473'
$ws.Range("E9").Value = '  -1.59%  '

$ws.Range("D10").Value = '7.45'
$ws.Range("E10").Value = '  -1.88%  '

$ws.Range("E11").Value = '  -5.45%  '

$ws.Range("D12").Value = '0.376'
$ws.Range("E12").Value = '  -2.04%  '

$ws.Range("D13").Value = '3.908.79'
$ws.Range("E13").Value = '  -2.49%  '

$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("D15").Value = '3.334.37'
$ws.Range("E15").Value = '  -2.51%  '

$ws.Range("D16").Value = '0.0000169'
$ws.Range("E16").Value = '  -4.84%  '

$ws.Range("D17").Value = '24.65'
$ws.Range("E17").Value = '  -1.33%  '

$ws.Range("D18").Value = '60.014.51'
$ws.Range("E18").Value = '  -5.84%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '5.68'
$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '13.43'
$ws.Range("E20").Value = '  +0.79%  '

$ws.Range("D21").Value = '9.03'
$ws.Range("E21").Value = '  -8.61%  '

$ws.Range("D22").Value = '354.41'
$ws.Range("E22").Value = '  -7.84%  '

$ws.Range("D23").Value = '0.558'
$ws.Range("E23").Value = '  -1.19%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = '3.472.33'
$ws.Range("E25").Value = '  -2.54%  '

$ws.Range("D26").Value = '69.24'
$ws.Range("E26").Value = '  -6.23%  '

$ws.Range("D27").Value = '0.0000112'
$ws.Range("E27").Value = '  +1.56%  '

$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '1.55'
$ws.Range("E28").Value = '  +9.14%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.46'
$ws.Range("E29").Value = '  +5.73%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D31").Value = '7.89'
$ws.Range("E31").Value = '  -1.07%  '

$ws.Range("D32").Value = '2.12'
$ws.Range("E32").Value = '  -4.07%  '

$ws.Range("E33").Value = '  -1.18%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").Value = '3.374.40'

$ws.Range("D36").Value = '22.95'
$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("D37").Value = '5.33'
$ws.Range("E37").Value = '  +2.79%  '

$ws.Range("D38").Value = '6.86'
$ws.Range("E38").Value = '  +1.17%  '

$ws.Range("D39").Value = '1.50'
$ws.Range("E39").Value = '  -0.40%  '

$ws.Range("D40").Value = '158.58'
$ws.Range("E40").Value = '  -3.25%  '

$ws.Range("D41").Value = '0.0761'
$ws.Range("E41").Value = '  -1.72%  '

$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.12%  '

$ws.Range("D43").Value = '4.39'
$ws.Range("E43").Value = '  +1.05%  '

$ws.Range("D44").Value = '1.20'
$ws.Range("E44").Value = '  +7.47%  '

$ws.Range("D45").Value = '40.86'

$ws.Range("D46").Value = '0.750'
$ws.Range("E46").Value = '  -4.94%  '

$ws.Range("D47").Value = '23.61'
$ws.Range("E47").Value = '  +0.53%  '

$ws.Range("D48").Value = '1.57'
$ws.Range("E48").Value = '  -2.94%  '

$ws.Range("D49").Value = '6.82'
$ws.Range("E49").Value = '  +1.42%  '

$ws.Range("D50").Value = '22.38'
$ws.Range("E50").Value = '  +9.95%  '

$ws.Range("D51").Value = '2.43'
$ws.Range("E51").Value = '  +16.98%  '

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
